$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.992.12"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").Value = "'1.857.85"
$ws.Range("E3").Value = "  -1.34%  "

$ws.Range("E4").Value = "  +0.39%  "

$ws.Range("D5").Value = "'312.13"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("D7").Value = "'0.5086"
$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("D8").Value = "'0.3831"
$ws.Range("E8").Value = "  -0.65%  "

$ws.Range("D9").Value = "'0.08231"
$ws.Range("E9").Value = "  -9.15%  "

$ws.Range("D10").Value = "'1.111"
$ws.Range("E10").Value = "  -1.55%  "

$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("D12").Value = "'6.204"
$ws.Range("E12").Value = "  -2.66%  "

$ws.Range("D13").Value = "'20.56"
$ws.Range("E13").Value = "  -1.35%  "

$ws.Range("D14").Value = "'1.861.79"
$ws.Range("E14").Value = "  -0.63%  "

$ws.Range("D15").Value = "'7.239"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("E16").Value = "  +0.38%  "

$ws.Range("D17").Value = "'0.00001096"
$ws.Range("E17").Value = "  -1.53%  "

$ws.Range("D18").Value = "'90.70"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").Value = "'0.06641"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").Value = "'17.66"
$ws.Range("E20").Value = "  -3.36%  "

$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").Value = "'6.014"
$ws.Range("E22").Value = "  -1.96%  "

$ws.Range("D23").Value = "'28.018.37"
$ws.Range("E23").Value = "  -0.54%  "

$ws.Range("D24").Value = "'11.06"

$ws.Range("D25").Value = "'2.233"
$ws.Range("E25").Value = "  -1.42%  "

$ws.Range("D26").Value = "'2.072.39"
$ws.Range("E26").Value = "  -0.83%  "

$ws.Range("D27").Value = "'2.508"
$ws.Range("E27").Value = "  -1.97%  "

$ws.Range("D28").Value = "'157.28"
$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("D29").Value = "'20.46"
$ws.Range("E29").Value = "  -2.00%  "

$ws.Range("D30").Value = "'124.75"
$ws.Range("E30").Value = "  -1.93%  "

$ws.Range("D31").Value = "'0.1058"
$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("D32").Value = "'1.032"
$ws.Range("E32").Value = "  -3.09%  "

$ws.Range("D33").Value = "'5.905"
$ws.Range("E33").Value = "  +4.83%  "

$ws.Range("D34").Value = "'3.592"
$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("D35").Value = "'9.378"
$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("D36").Value = "'0.06509"
$ws.Range("E36").Value = "  -1.76%  "

$ws.Range("D37").Value = "'0.02411"
$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("D38").Value = "'0.2172"
$ws.Range("E38").Value = "  -1.18%  "

$ws.Range("D39").Value = "'0.6544"
$ws.Range("E39").Value = "  +1.55%  "

$ws.Range("E40").Value = "  -1.56%  "

$ws.Range("D41").Value = "'4.985"
$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("D42").Value = "'1.219"
$ws.Range("E42").Value = "  -5.71%  "

$ws.Range("D43").Value = "'11.14"
$ws.Range("E43").Value = "  -3.39%  "

$ws.Range("D44").Value = "'0.6126"
$ws.Range("E44").Value = "  +1.15%  "

$ws.Range("D45").Value = "'13.14"
$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("D46").Value = "'1.283"
$ws.Range("E46").Value = "  +0.63%  "

$ws.Range("D47").Value = "'3.651"
$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("D48").Value = "'2.012"
$ws.Range("E48").Value = "  +0.14%  "

$ws.Range("D49").Value = "'1.206"
$ws.Range("E49").Value = "  -2.98%  "

$ws.Range("D50").Value = "'119.95"
$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("D51").Value = "'78.29"
$ws.Range("E51").Value = "  -1.89%  "
